# Updated AVGO and MRVL data to the latest quarter.
#
# - Relabel the quarter headers from "Qn '2x" to "Qn FYxx".
# - Add two new quarterly columns (G = Q2 FY26, H = Q3 FY26).
# - Fill in report date for Q1 FY25 (col E) which was previously blank,
#   and add report dates for the two new quarters.
# - Add the new quarter's KPI figures for every product line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: quarter labels ---
$ws.Range("B1").Value = "Q1 FY25"
$ws.Range("C1").Value = "Q2 FY25"
$ws.Range("D1").Value = "Q3 FY25"
$ws.Range("E1").Value = "Q4 FY25"
$ws.Range("F1").Value = "Q1 FY26"
$ws.Range("G1").Value = "Q2 FY26"
$ws.Range("H1").Value = "Q3 FY26"

# --- Report Date row: carry the existing date format (col D) onto the
#     newly populated cells so the same built-in date numFmt is reused. ---
$ws.Range("D2").Copy()
$ws.Range("F2:H2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E2").Value = 45689
$ws.Range("F2").Value = 45780
$ws.Range("G2").Value = 45871
$ws.Range("H2").Value = 45962

# --- KPI rows: carry the existing "0.0" numeric format (col F) onto the
#     new columns before writing values. Note G5 is intentionally left
#     with the default (General) number format, matching the source data,
#     so its format is not copied there. ---
$ws.Range("F3:F4").Copy()
$ws.Range("G3:H4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F5").Copy()
$ws.Range("H5").PasteSpecial(-4122)     # xlPasteFormats

$ws.Range("F6:F7").Copy()
$ws.Range("G6:H7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G3").Value = 1490.5
$ws.Range("H3").Value = 1517.9

$ws.Range("G4").Value = 193.6
$ws.Range("H4").Value = 237.2

$ws.Range("G5").Value = 130.1
$ws.Range("H5").Value = 167.8

$ws.Range("G6").Value = 115.9
$ws.Range("H6").Value = 116.6

$ws.Range("G7").Value = 76
$ws.Range("H7").Value = 35

$excel.CutCopyMode = $false

# Selection mirrors the authored workbook (row 3 selected, A3:H3).
$ws.Range("A3:H3").Select() | Out-Null

